# The deck currently carries the "Integral" design (theme1.xml, used by the
# slide master / notes master / handout master) and an unused "Office Theme"
# palette (theme2.xml). The author swapped the two themes' content so the
# deck now uses the standard "Office" colour scheme while "Integral" is set
# aside -- i.e. theme1.xml's <a:clrScheme> becomes the Office palette and
# theme2.xml's becomes the Integral palette.
#
# fontScheme/fmtScheme are identical between the two themes, so the only
# real content delta is the 12 colour slots of <a:clrScheme>. We recolour
# the presentation's theme (reachable through SlideMaster.Theme) to the
# "Office" palette using ThemeColorScheme, PowerPoint's supported API for
# editing theme colours (Master.Theme.ThemeColorScheme.Colors(i).RGB).

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Office theme colour scheme, in Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Colors(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1      000000
$tcs.Colors(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Colors(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Colors(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Colors(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Colors(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Colors(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Colors(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Colors(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Colors(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink 954F72
